$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = '30.676.26'
$ws.Range("E2").Value2 = '  +1.96%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = '1.894.46'
$ws.Range("E3").Value2 = '  +1.14%  '

# Row 4
$ws.Range("E4").Value2 = '  +0.15%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '241.62'
$ws.Range("E5").Value2 = '  -0.13%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '1.001'
$ws.Range("E6").Value2 = '  +0.09%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = '0.4921'
$ws.Range("E7").Value2 = '  +1.23%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = '0.2933'
$ws.Range("E8").Value2 = '  +1.80%  '

# Row 9
$ws.Range("E9").Value2 = '  +2.60%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '1.895.39'
$ws.Range("E10").Value2 = '  +1.22%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '17.16'
$ws.Range("E11").Value2 = '  +5.48%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = '0.07256'
$ws.Range("E12").Value2 = '  +1.04%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = '90.71'
$ws.Range("E13").Value2 = '  +5.85%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '0.6747'
$ws.Range("E14").Value2 = '  +2.08%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '5.020'
$ws.Range("E15").Value2 = '  +2.51%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '30.673.10'
$ws.Range("E16").Value2 = '  +2.11%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = '0.000007947'
$ws.Range("E17").Value2 = '  +2.77%  '

# Row 18
$ws.Range("E18").Value2 = '  +0.09%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '13.08'
$ws.Range("E19").Value2 = '  +3.15%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '2.139.84'

# Row 21
$ws.Range("E21").Value2 = '  +0.02%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '4.795'
$ws.Range("E22").Value2 = '  +1.43%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '188.85'
$ws.Range("E23").Value2 = '  +32.78%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '6.068'
$ws.Range("E24").Value2 = '  +4.32%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '9.348'
$ws.Range("E25").Value2 = '  +2.46%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '156.76'
$ws.Range("E26").Value2 = '  +3.48%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = '18.74'
$ws.Range("E27").Value2 = '  +11.09%  '

# Row 28
$ws.Range("E28").Value2 = '  +1.28%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '1.401'
$ws.Range("E29").Value2 = '  +0.88%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '4.273'
$ws.Range("E30").Value2 = '  +2.36%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = '0.09086'
$ws.Range("E31").Value2 = '  +3.83%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '3.989'
$ws.Range("E32").Value2 = '  +0.50%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '0.05230'
$ws.Range("E33").Value2 = '  +2.78%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '0.7364'
$ws.Range("E34").Value2 = '  +4.00%  '

# Row 35
$ws.Range("E35").Value2 = '  +0.34%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '2.752'
$ws.Range("E36").Value2 = '  +3.06%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '0.01830'
$ws.Range("E37").Value2 = '  +0.10%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '2.678'
$ws.Range("E38").Value2 = '  +0.18%  '

# Row 39
$ws.Range("B39").Value2 = 'TrustWalletToken'
$ws.Range("C39").Value2 = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '0.9316'
$ws.Range("E39").Value2 = '  +0.98%  '

# Row 40
$ws.Range("B40").Value2 = 'RenderToken'
$ws.Range("C40").Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '2.120'
$ws.Range("E40").Value2 = '  -0.69%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '0.4384'
$ws.Range("E41").Value2 = '  +4.39%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '104.95'
$ws.Range("E42").Value2 = '  +1.33%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '1.001'
$ws.Range("E43").Value2 = '  +0.19%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '5.723'
$ws.Range("E44").Value2 = '  -0.16%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '7.515'
$ws.Range("E45").Value2 = '  +1.95%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '0.1346'
$ws.Range("E46").Value2 = '  +5.83%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '0.05859'
$ws.Range("E47").Value2 = '  +2.61%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = '8.725'
$ws.Range("E48").Value2 = '  +5.77%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '1.421'
$ws.Range("E49").Value2 = '  +6.90%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '0.3927'
$ws.Range("E50").Value2 = '  +5.32%  '

# Row 51
$ws.Range("E51").Value2 = '  +3.16%  '
